$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Paper", "AMD",  100, "2025-08-11 09:15", 171.35, "2025-08-11 10:05", 174.78, 343,  2,     "Entry 171.35, Target 174.78, Stop 169.64"),
    @("Paper", "ETHU", 100, "2025-08-12 11:28", 162.78, "2025-08-12 12:00", 166.04, 326,  2,     "Entry 162.78, Target 166.04, Stop 161.15"),
    @("Paper", "AMD",  100, "2025-08-18 09:55", 72.49,  "2025-08-18 11:05", 71.77,  -72,  -0.99, "LABU Entry 72.49, Target 73.94, Stop 71.77")
)

$startRow = 21
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
}
